# Apply NGBD bond type/amount rows to Sheet1 ("Case_Data")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format first for cells whose values would otherwise be
# auto-coerced to numbers by Excel (e.g. "4510.11", "50", "$ 0").
$ws.Range("D830:D831").NumberFormat = "@"
$ws.Range("H834:K836").NumberFormat = "@"

$ws.Range("A819").Value = "22CRB00136"
$ws.Range("B819").Value = "Hemmeter"
$ws.Range("C819").Value = "DOMESTIC VIOLENCE"
$ws.Range("D819").Value = "2919.25(A)"
$ws.Range("E819").Value = "No Data"
$ws.Range("F819").Value = "Not Guilty"

$ws.Range("A820").Value = "22CRB00136"
$ws.Range("B820").Value = "Hemmeter"
$ws.Range("C820").Value = "ASSAULT - M1"
$ws.Range("D820").Value = "2903.13(A)"
$ws.Range("E820").Value = "No Data"
$ws.Range("F820").Value = "Not Guilty"

$ws.Range("A821").Value = "22CRB00136"
$ws.Range("B821").Value = "Hemmeter"
$ws.Range("C821").Value = "DOMESTIC VIOLENCE"
$ws.Range("D821").Value = "2919.25(A)"
$ws.Range("E821").Value = "No Data"
$ws.Range("F821").Value = "Not Guilty"

$ws.Range("A822").Value = "22CRB00136"
$ws.Range("B822").Value = "Hemmeter"
$ws.Range("C822").Value = "ASSAULT - M1"
$ws.Range("D822").Value = "2903.13(A)"
$ws.Range("E822").Value = "No Data"
$ws.Range("F822").Value = "Not Guilty"

$ws.Range("A823").Value = "22CRB00136"
$ws.Range("B823").Value = "Hemmeter"
$ws.Range("C823").Value = "DOMESTIC VIOLENCE"
$ws.Range("D823").Value = "2919.25(A)"
$ws.Range("E823").Value = "No Data"
$ws.Range("F823").Value = "Not Guilty"

$ws.Range("A824").Value = "22CRB00136"
$ws.Range("B824").Value = "Hemmeter"
$ws.Range("C824").Value = "ASSAULT - M1"
$ws.Range("D824").Value = "2903.13(A)"
$ws.Range("E824").Value = "No Data"
$ws.Range("F824").Value = "Not Guilty"

$ws.Range("A825").Value = "22CRB00142"
$ws.Range("B825").Value = "Hemmeter"
$ws.Range("C825").Value = "THEFT / M1"
$ws.Range("D825").Value = "2913.02(A)(1)*"
$ws.Range("E825").Value = "M1"
$ws.Range("F825").Value = "Not Guilty"

$ws.Range("A826").Value = "22CRB00136"
$ws.Range("B826").Value = "Hemmeter"
$ws.Range("C826").Value = "DOMESTIC VIOLENCE"
$ws.Range("D826").Value = "2919.25(A)"
$ws.Range("E826").Value = "No Data"
$ws.Range("F826").Value = "Not Guilty"

$ws.Range("A827").Value = "22CRB00136"
$ws.Range("B827").Value = "Hemmeter"
$ws.Range("C827").Value = "ASSAULT - M1"
$ws.Range("D827").Value = "2903.13(A)"
$ws.Range("E827").Value = "No Data"
$ws.Range("F827").Value = "Not Guilty"

$ws.Range("A828").Value = "22CRB00142"
$ws.Range("B828").Value = "Hemmeter"
$ws.Range("C828").Value = "THEFT / M1"
$ws.Range("D828").Value = "2913.02(A)(1)*"
$ws.Range("E828").Value = "M1"
$ws.Range("F828").Value = "Not Guilty"

$ws.Range("A829").Value = "22CRB00142"
$ws.Range("B829").Value = "Hemmeter"
$ws.Range("C829").Value = "THEFT / M1"
$ws.Range("D829").Value = "2913.02(A)(1)*"
$ws.Range("E829").Value = "M1"
$ws.Range("F829").Value = "Not Guilty"

$ws.Range("B830").Value = "Hemmeter"
$ws.Range("C830").Value = "Driving Under Suspension"
$ws.Range("D830").Value = "4510.11"
$ws.Range("E830").Value = "M1"
$ws.Range("F830").Value = "Not Guilty"

$ws.Range("B831").Value = "Hemmeter"
$ws.Range("C831").Value = "Driving Under Suspension"
$ws.Range("D831").Value = "4510.11"
$ws.Range("E831").Value = "M1"
$ws.Range("F831").Value = "Not Guilty"

$ws.Range("A832").Value = "22CRB00136"
$ws.Range("B832").Value = "Hemmeter"
$ws.Range("C832").Value = "DOMESTIC VIOLENCE"
$ws.Range("D832").Value = "2919.25(A)"
$ws.Range("E832").Value = "No Data"
$ws.Range("F832").Value = "Not Guilty"

$ws.Range("A833").Value = "22CRB00136"
$ws.Range("B833").Value = "Hemmeter"
$ws.Range("C833").Value = "ASSAULT - M1"
$ws.Range("D833").Value = "2903.13(A)"
$ws.Range("E833").Value = "No Data"
$ws.Range("F833").Value = "Not Guilty"

$ws.Range("A834").Value = "21CRB00626"
$ws.Range("B834").Value = "Hemmeter"
$ws.Range("C834").Value = "CRIMINAL MISCHIEF"
$ws.Range("D834").Value = "2909.07(A)(1)"
$ws.Range("E834").Value = "M3"
$ws.Range("F834").Value = "No Contest"
$ws.Range("G834").Value = "Guilty"
$ws.Range("H834").Value = "`$ 0"
$ws.Range("I834").Value = "`$ 0"
$ws.Range("J834").Value = "50"
$ws.Range("K834").Value = "None"

$ws.Range("A835").Value = "21CRB00626"
$ws.Range("B835").Value = "Hemmeter"
$ws.Range("C835").Value = "ASSAULT - M1"
$ws.Range("D835").Value = "2903.13(A)"
$ws.Range("E835").Value = "M1"
$ws.Range("F835").Value = "No Contest"
$ws.Range("G835").Value = "Guilty"
$ws.Range("H835").Value = "`$ 0"
$ws.Range("I835").Value = "`$ 0"
$ws.Range("J835").Value = "None"
$ws.Range("K835").Value = "None"

$ws.Range("A836").Value = "21CRB00626"
$ws.Range("B836").Value = "Hemmeter"
$ws.Range("C836").Value = "DISORDERLY CONDUCT"
$ws.Range("D836").Value = "2917.11A1"
$ws.Range("E836").Value = "MM"
$ws.Range("F836").Value = "No Contest"
$ws.Range("G836").Value = "Guilty"
$ws.Range("H836").Value = "`$ 0"
$ws.Range("I836").Value = "`$ 0"
$ws.Range("J836").Value = "None"
$ws.Range("K836").Value = "None"

Write-Host "Added rows 819-836 to Sheet1"
